# fix(publipostage): Try to solve Excel emoji problem
#
# The "statut" column (column A) stores a small set of emoji markers that
# are repeated across every data row. Replace each emoji with its plain
# text equivalent everywhere it occurs:
#   📕 -> -3
#   📘 -> ⚠️
#   📗 -> ✅
#   📙 -> +3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "📕" = "-3"
    "📘" = "⚠️"
    "📗" = "✅"
    "📙" = "+3"
}
# Values that look numeric need an explicit text format so Excel keeps
# storing them as text (matching the original string-typed cells) instead
# of silently converting them to numbers.
$numericLooking = @("-3", "+3")

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Find the "statut" column from the header row instead of hardcoding it.
$firstCol = $used.Column
$lastCol = $firstCol + $used.Columns.Count - 1
$statutCol = $firstCol
for ($c = $firstCol; $c -le $lastCol; $c++) {
    if ($ws.Cells.Item($firstRow, $c).Value2 -eq "statut") {
        $statutCol = $c
        break
    }
}

for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $statutCol)
    $val = $cell.Value2
    if ($map.ContainsKey($val)) {
        $newVal = $map[$val]
        if ($numericLooking -contains $newVal) {
            $cell.NumberFormat = "@"
        }
        $cell.Value2 = $newVal
    }
}
